# Adding new progress as of date 04 Nov 2025.
# For the "Training Dashboard" sheet, rows 3-29:
#   - PERIOD TO EXPIRE (column H) decreases by 1 day.
#   - LAST UPDATE (column I) moves from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 29; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE
    # Assign through a text formula and then collapse it back down to a
    # plain value via Copy/PasteSpecial so the cell keeps storing a literal
    # text string ("04-Nov-2025") instead of being auto-converted into a
    # real date serial number (which would also change its number format).
    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy()
    $iCell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
